$wb = $excel.ActiveWorkbook

# Sheet3: change A3 from "Cast iron" to "Bronze", and move the selection to D10
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A3").Value = "Bronze"
$ws3.Range("D10").Select()

# Sheet2: move the selection to D13 (no data change)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("D13").Select()

# Leave Sheet3 as the active/selected sheet (matches tabSelected in the diff)
$ws3.Activate()
$ws3.Range("D10").Select()
